# Insert two new price records right before the current row 372
# (Fecha = 2022-03-17 / serial 44637), pushing the existing rows
# 372..442 down to 374..444, matching the new data snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 372 and onward down by inserting 2 blank rows at 372.
$ws.Rows.Item(372).Resize(2).Insert()

# --- New row 372: "1a plateado" ---
$ws.Cells.Item(372, 1).Value = 4
$ws.Cells.Item(372, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(372, 3).Value = "Los Lagos"
$ws.Cells.Item(372, 4).Value = 44637
$ws.Cells.Item(372, 5).Value = 10
$ws.Cells.Item(372, 6).Value = "Fruta"
$ws.Cells.Item(372, 7).Value = 100102
$ws.Cells.Item(372, 8).Value = "Cítricos"
$ws.Cells.Item(372, 9).Value = 100102003
$ws.Cells.Item(372, 10).Value = "Limón"
$ws.Cells.Item(372, 11).Value = "Sin especificar"
$ws.Cells.Item(372, 12).Value = "1a plateado"
$ws.Cells.Item(372, 13).Value = 800
$ws.Cells.Item(372, 14).Value = 27000
$ws.Cells.Item(372, 15).Value = 28000
$ws.Cells.Item(372, 16).Value = 27500
$ws.Cells.Item(372, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(372, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(372, 19).Value = 1528
$ws.Cells.Item(372, 20).Value = 18

# --- New row 373: "2a plateado" ---
$ws.Cells.Item(373, 1).Value = 4
$ws.Cells.Item(373, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(373, 3).Value = "Los Lagos"
$ws.Cells.Item(373, 4).Value = 44637
$ws.Cells.Item(373, 5).Value = 10
$ws.Cells.Item(373, 6).Value = "Fruta"
$ws.Cells.Item(373, 7).Value = 100102
$ws.Cells.Item(373, 8).Value = "Cítricos"
$ws.Cells.Item(373, 9).Value = 100102003
$ws.Cells.Item(373, 10).Value = "Limón"
$ws.Cells.Item(373, 11).Value = "Sin especificar"
$ws.Cells.Item(373, 12).Value = "2a plateado"
$ws.Cells.Item(373, 13).Value = 300
$ws.Cells.Item(373, 14).Value = 25000
$ws.Cells.Item(373, 15).Value = 25000
$ws.Cells.Item(373, 16).Value = 25000
$ws.Cells.Item(373, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(373, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(373, 19).Value = 1389
$ws.Cells.Item(373, 20).Value = 18

# Match the date-column style used elsewhere (column D, rows 4..442) for
# the two new rows.
$ws.Range("D372:D373").NumberFormat = $ws.Range("D374").NumberFormat
